$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.2708946666666667
$ws.Cells.Item(2, 8).Value = 0.812684
$ws.Cells.Item(2, 9).Value = 0.1616296696421007
$ws.Cells.Item(2, 10).Value = 0.1616296696421007
$ws.Cells.Item(2, 15).Value = 0.9347132976570145
$ws.Cells.Item(2, 16).Value = 0.9347132976570145
$ws.Cells.Item(2, 17).Value = 2.406284272738222
$ws.Cells.Item(2, 18).Value = 21.656558454644
$ws.Cells.Item(2, 19).Value = 0.1510774015103817
$ws.Cells.Item(2, 20).Value = 0.1510774015103817
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.2708946666666667
$ws.Cells.Item(3, 8).Value = 0.812684
$ws.Cells.Item(3, 9).Value = 0.1616296696421007
$ws.Cells.Item(3, 10).Value = 0.1616296696421007
$ws.Cells.Item(3, 13).Value = 0.616144
$ws.Cells.Item(3, 14).Value = 1.848432
$ws.Cells.Item(3, 15).Value = 0.06483569448352988
$ws.Cells.Item(3, 16).Value = 0.0648356944835299
$ws.Cells.Item(3, 17).Value = 0.1669101234986667
$ws.Cells.Item(3, 18).Value = 1.502191111488
$ws.Cells.Item(3, 19).Value = 0.0104793718803891
$ws.Cells.Item(3, 20).Value = 0.0104793718803891
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.2708946666666667
$ws.Cells.Item(4, 8).Value = 0.812684
$ws.Cells.Item(4, 9).Value = 0.1616296696421007
$ws.Cells.Item(4, 10).Value = 0.1616296696421007
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.004286
$ws.Cells.Item(4, 14).Value = 0.012858
$ws.Cells.Item(4, 15).Value = 0.0004510078594555965
$ws.Cells.Item(4, 16).Value = 0.0004510078594555965
$ws.Cells.Item(4, 17).Value = 0.001161054541333333
$ws.Cells.Item(4, 18).Value = 0.010449490872
$ws.Cells.Item(4, 19).Value = 0.00007289625132979903
$ws.Cells.Item(4, 20).Value = 0.00007289625132979903
$ws.Cells.Item(5, 9).Value = 0.6313295261673385
$ws.Cells.Item(5, 10).Value = 0.6313295261673384
$ws.Cells.Item(5, 15).Value = 0.9347132976570145
$ws.Cells.Item(5, 16).Value = 0.9347132976570145
$ws.Cells.Item(5, 17).Value = 9.399006463947112
$ws.Cells.Item(5, 19).Value = 0.5901121033121134
$ws.Cells.Item(5, 20).Value = 0.5901121033121133
$ws.Cells.Item(6, 9).Value = 0.6313295261673385
$ws.Cells.Item(6, 10).Value = 0.6313295261673384
$ws.Cells.Item(6, 13).Value = 0.616144
$ws.Cells.Item(6, 14).Value = 1.848432
$ws.Cells.Item(6, 15).Value = 0.06483569448352988
$ws.Cells.Item(6, 16).Value = 0.0648356944835299
$ws.Cells.Item(6, 17).Value = 0.6519551108053334
$ws.Cells.Item(6, 18).Value = 5.867595997248
$ws.Cells.Item(6, 19).Value = 0.04093268827701724
$ws.Cells.Item(6, 20).Value = 0.04093268827701724
$ws.Cells.Item(7, 9).Value = 0.6313295261673385
$ws.Cells.Item(7, 10).Value = 0.6313295261673384
$ws.Cells.Item(7, 11).Value = 2
$ws.Cells.Item(7, 12).Value = 0.6666666666666666
$ws.Cells.Item(7, 13).Value = 0.004286
$ws.Cells.Item(7, 14).Value = 0.012858
$ws.Cells.Item(7, 15).Value = 0.0004510078594555965
$ws.Cells.Item(7, 16).Value = 0.0004510078594555965
$ws.Cells.Item(7, 17).Value = 0.004535108034666667
$ws.Cells.Item(7, 18).Value = 0.040815972312
$ws.Cells.Item(7, 19).Value = 0.0002847345782078474
$ws.Cells.Item(7, 20).Value = 0.0002847345782078473
$ws.Cells.Item(8, 9).Value = 0.2070408041905609
$ws.Cells.Item(8, 10).Value = 0.2070408041905609
$ws.Cells.Item(8, 15).Value = 0.9347132976570145
$ws.Cells.Item(8, 16).Value = 0.9347132976570145
$ws.Cells.Item(8, 19).Value = 0.1935237928345194
$ws.Cells.Item(8, 20).Value = 0.1935237928345194
$ws.Cells.Item(9, 9).Value = 0.2070408041905609
$ws.Cells.Item(9, 10).Value = 0.2070408041905609
$ws.Cells.Item(9, 13).Value = 0.616144
$ws.Cells.Item(9, 14).Value = 1.848432
$ws.Cells.Item(9, 15).Value = 0.06483569448352988
$ws.Cells.Item(9, 16).Value = 0.0648356944835299
$ws.Cells.Item(9, 17).Value = 0.2138048433386667
$ws.Cells.Item(9, 18).Value = 1.924243590048
$ws.Cells.Item(9, 19).Value = 0.01342363432612354
$ws.Cells.Item(9, 20).Value = 0.01342363432612354
$ws.Cells.Item(10, 9).Value = 0.2070408041905609
$ws.Cells.Item(10, 10).Value = 0.2070408041905609
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.004286
$ws.Cells.Item(10, 14).Value = 0.012858
$ws.Cells.Item(10, 15).Value = 0.0004510078594555965
$ws.Cells.Item(10, 16).Value = 0.0004510078594555965
$ws.Cells.Item(10, 17).Value = 0.001487262001333333
$ws.Cells.Item(10, 18).Value = 0.013385358012
$ws.Cells.Item(10, 19).Value = 0.00009337702991795017
$ws.Cells.Item(10, 20).Value = 0.00009337702991795017
